# Updated cryptos list on Sat Jul 22 02:55:55 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.960.00"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.894.39"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - XRP
$ws.Range("D5").Value = "'0.7734"
$ws.Range("E5").Value = "  -3.09%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'244.71"

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.00%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3139"
$ws.Range("E8").Value = "  -1.28%  "

# Row 9 - Solana
$ws.Range("E9").Value = "  +0.97%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07266"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.08373"
$ws.Range("E11").Value = "  +3.16%  "

# Row 12 - Polygon
$ws.Range("D12").Value = "'0.7738"
$ws.Range("E12").Value = "  +0.18%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'5.479"
$ws.Range("E13").Value = "  -2.16%  "

# Row 14 - was WrappedEther, now Litecoin
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'95.07"
$ws.Range("E14").Value = "  +2.26%  "

# Row 15 - was Litecoin, now WrappedEther
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.824.68"
$ws.Range("E15").Value = "  -3.50%  "

# Row 16 - Uniswap
$ws.Range("D16").Value = "'6.200"
$ws.Range("E16").Value = "  +0.23%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "29.972.08"
$ws.Range("E17").Value = "  +0.23%  "

# Row 18 - Avalanche
$ws.Range("E18").Value = "  +0.17%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'247.30"
$ws.Range("E19").Value = "  +0.59%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.000007875"
$ws.Range("E20").Value = "  +1.32%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'8.143"
$ws.Range("E21").Value = "  -1.88%  "

# Row 22 - was Dai, now WrappedliquidstakedEther2.0
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.142.27"
$ws.Range("E22").Value = "  -0.43%  "

# Row 23 - was WrappedliquidstakedEther2.0, now Dai
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.0000"
$ws.Range("E23").Value = "  -0.04%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.02%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "'0.1598"
$ws.Range("E25").Value = "  -4.67%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'9.548"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'162.56"
$ws.Range("E27").Value = "  -1.08%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'18.79"
$ws.Range("E28").Value = "  +0.13%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'2.040"
$ws.Range("E29").Value = "  -1.68%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +0.70%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.10%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.559"
$ws.Range("E32").Value = "  +1.06%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'4.109"
$ws.Range("E33").Value = "  +0.51%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "'0.05475"
$ws.Range("E34").Value = "  -2.91%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  -3.17%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.7523"
$ws.Range("E36").Value = "  +0.97%  "

# Row 37 - Frax
$ws.Range("D37").Value = "'1.003"
$ws.Range("E37").Value = "  +0.03%  "

# Row 38 - HuobiToken
$ws.Range("D38").Value = "'2.673"
$ws.Range("E38").Value = "  +1.38%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.14%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "'2.791"
$ws.Range("E40").Value = "  +0.08%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.4496"
$ws.Range("E41").Value = "  +1.24%  "

# Row 42 - Aave
$ws.Range("D42").Value = "'74.18"
$ws.Range("E42").Value = "  -0.95%  "

# Row 43 - was Maker, now FraxShare
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.085"
$ws.Range("E43").Value = "  +2.11%  "

# Row 44 - was FraxShare, now Maker
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.095.75"
$ws.Range("E44").Value = "  -6.74%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "'0.8529"
$ws.Range("E45").Value = "  -0.35%  "

# Row 46 - PaxDollar
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +0.29%  "

# Row 48 - Quant
$ws.Range("D48").Value = "'102.73"
$ws.Range("E48").Value = "  -1.87%  "

# Row 49 - Aptos
$ws.Range("D49").Value = "'7.608"
$ws.Range("E49").Value = "  +1.39%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'9.825"
$ws.Range("E50").Value = "  -3.05%  "

# Row 51 - SynthetixNetwork
$ws.Range("D51").Value = "'3.014"
$ws.Range("E51").Value = "  +1.12%  "
